$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 247 (Femacal de La Calera / Berenjena block),
# pushing the existing rows 247-357 down to 248-358.
$ws.Rows.Item(247).Insert()

# Populate the newly inserted row 247 with the new weekly data point.
$ws.Range("A247").Value = 3
$ws.Range("B247").Value = "Femacal de La Calera"
$ws.Range("C247").Value = "Coquimbo"
$ws.Range("D247").Value = 44825
$ws.Range("E247").Value = 5
$ws.Range("F247").Value = 100112001
$ws.Range("G247").Value = "Berenjena"
$ws.Range("H247").Value = "Sin especificar"
$ws.Range("I247").Value = "Primera"
$ws.Range("J247").Value = 45
$ws.Range("K247").Value = 13000
$ws.Range("L247").Value = 13000
$ws.Range("M247").Value = 13000
$ws.Range("N247").Value = "$/caja 60 unidades"
$ws.Range("O247").Value = "Región de Arica y Parinacota"
$ws.Range("P247").Value = 217
$ws.Range("Q247").Value = 60
$ws.Range("R247").Value = "Hortaliza"
